$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unneeded trailing column (E) and trailing rows (8:11)
# first, so the remaining ranges line up with the final A1:D7 layout.
$ws.Range("E:E").Delete()
$ws.Range("8:11").Delete()

# Update header row (B1:D1 already carry the bold/border/center style,
# so only the text needs to change).
$ws.Range("B1").Value = "id"
$ws.Range("C1").Value = "raw_grade"
$ws.Range("D1").Value = "grade"

# Replace the date-based id column with plain integers 0..5.
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5

# raw_grade numeric codes 1..6
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 3
$ws.Range("B5").Value = 4
$ws.Range("B6").Value = 5
$ws.Range("B7").Value = 6

# grade letters
$ws.Range("C2").Value = "a"
$ws.Range("C3").Value = "b"
$ws.Range("C4").Value = "b"
$ws.Range("C5").Value = "a"
$ws.Range("C6").Value = "a"
$ws.Range("C7").Value = "e"

# grade descriptions
$ws.Range("D2").Value = "very good"
$ws.Range("D3").Value = "good"
$ws.Range("D4").Value = "good"
$ws.Range("D5").Value = "very good"
$ws.Range("D6").Value = "very good"
$ws.Range("D7").Value = "very bad"

# Column A used to carry a custom date-time number format (style index 2).
# The new id column should look like the header style instead (bold,
# bordered, centered, general number format) -- copy B1's formatting
# (which already is that exact style) onto A2:A7 so no new style / number
# format entries get created and the old date format becomes unused.
$ws.Range("B1").Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
